$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly price rows (row 2 and row 3) had their date/volume/price data swapped.
# Row 2 should now hold what was previously row 3's data (and vice versa) for
# columns D, J, K, L, M, P.

$ws.Range("D2").Value = 44749
$ws.Range("J2").Value = 90
$ws.Range("K2").Value = 17000
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = 17556
$ws.Range("P2").Value = 1170

$ws.Range("D3").Value = 44839
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 15000
$ws.Range("L3").Value = 16000
$ws.Range("M3").Value = 15600
$ws.Range("P3").Value = 1040
